$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.534.23"
$ws.Range("E2").Value = "  +3.11%  "
$ws.Range("D3").Value = "3.773.37"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.53"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.16"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("D7").Value = "3.771.36"
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.540"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("E10").Value = "  +4.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.42"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.71%  "
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.64"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.49%  "
$ws.Range("E14").Value = "  +4.27%  "
$ws.Range("D15").Value = "4.399.75"
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").Value = "3.769.08"
$ws.Range("E16").Value = "  +1.85%  "
$ws.Range("D17").Value = "69.506.85"
$ws.Range("E17").Value = "  +2.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.34"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.23%  "
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.20"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.93"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +17.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "498.09"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.733"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000155"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +11.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.62"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.36"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.46"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.36"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.47%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.01"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.52"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +7.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.01"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.14"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.91%  "
$ws.Range("D34").Value = "3.916.47"
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("B35").Value = "RenzoRestakedETH"
$ws.Range("C35").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D35").Value = "3.704.39"
$ws.Range("E35").Value = "  +1.57%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.109"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.41%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  +2.34%  "
$ws.Range("E39").Value = "  +3.20%  "
$ws.Range("E40").Value = "  +2.15%  "
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.08"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +10.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "443.42"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.34%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.01"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.63%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.72"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.52"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.48%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.84"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").Value = "2.822.31"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "140.95"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("E51").Value = "  +2.62%  "
